$d = $word.ActiveDocument

# --- Change 1: first paragraph (the **ID__...__ID** marker) ---
# Add paragraph border spacing (pBdr top/left/bottom/right w:space="5"),
# change left indent from 120 -> 225 twips (6pt -> 11.25pt), and update/merge
# the marker text into a single run with no trailing space.
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.Borders.DistanceFromTop = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromBottom = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromLeft = 5
$p1.Range.ParagraphFormat.Borders.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25

$d.Content.Find.Execute("**ID__AFFARS_pgi_5301_topic_17__ID** ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "**ID__AFFARS_AFICC_PGI_5301_601_91__ID**", 2)

# --- Change 2: table row height ---
# The "Competitive Awards" / "1-500" row grows from 1005 -> 1155 twips
# (50.25pt -> 57.75pt). The first column has a vertical merge, so index the
# row via a specific cell rather than Tables(1).Rows(2) directly.
$t = $d.Tables(1)
$cell = $t.Cell(2, 1)
$cell.Row.HeightRule = 1
$cell.Row.Height = 57.75

Write-Output "edit applied"
